$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.575558
$ws.Range("H2").Value = 10.726674
$ws.Range("I2").Value = 0.025194653521236
$ws.Range("J2").Value = 0.02519465352123599
$ws.Range("M2").Value = 0.06624833333333334
$ws.Range("N2").Value = 0.198745
$ws.Range("Q2").Value = 0.2368747582366666
$ws.Range("R2").Value = 2.13187282413
$ws.Range("S2").Value = 0.025194653521236
$ws.Range("T2").Value = 0.02519465352123599

# Row 3
$ws.Range("I3").Value = 0.7460690747908298
$ws.Range("J3").Value = 0.7460690747908298
$ws.Range("M3").Value = 0.06624833333333334
$ws.Range("N3").Value = 0.198745
$ws.Range("Q3").Value = 7.014382300196112
$ws.Range("R3").Value = 63.12944070176501
$ws.Range("S3").Value = 0.7460690747908298
$ws.Range("T3").Value = 0.7460690747908298

# Row 4
$ws.Range("G4").Value = 32.36130266666667
$ws.Range("H4").Value = 97.08390800000001
$ws.Range("I4").Value = 0.2280292497513723
$ws.Range("J4").Value = 0.2280292497513723
$ws.Range("M4").Value = 0.06624833333333334
$ws.Range("N4").Value = 0.198745
$ws.Range("Q4").Value = 2.143882366162222
$ws.Range("R4").Value = 19.29494129546
$ws.Range("S4").Value = 0.2280292497513723
$ws.Range("T4").Value = 0.2280292497513723

# Row 5
$ws.Range("G5").Value = 0.1003386666666667
$ws.Range("H5").Value = 0.301016
$ws.Range("I5").Value = 0.000707021936561918
$ws.Range("J5").Value = 0.0007070219365619179
$ws.Range("M5").Value = 0.06624833333333334
$ws.Range("N5").Value = 0.198745
$ws.Range("Q5").Value = 0.006647269435555557
$ws.Range("R5").Value = 0.05982542492
$ws.Range("S5").Value = 0.000707021936561918
$ws.Range("T5").Value = 0.0007070219365619179
